$wb = $excel.ActiveWorkbook

# --- SignIn sheet: add two new rows (hyperlinked emails + passwords) ---
$signIn = $wb.Worksheets.Item("SignIn")

$signIn.Range("A6").Value = "bmacedo1987@gmail.com"
$signIn.Range("B6").Value = "industryconnect"
$signIn.Hyperlinks.Add($signIn.Range("A6"), "mailto:bmacedo1987@gmail.com") | Out-Null

$signIn.Range("A7").Value = "johndoe@gmail.com"
$signIn.Range("B7").Value = 123456
$signIn.Hyperlinks.Add($signIn.Range("A7"), "mailto:johndoe@gmail.com") | Out-Null

# --- Update selection / active sheet state ---
$signUp = $wb.Worksheets.Item("SignUp")
$signUp.Range("A6").Select()

$signIn.Activate()
$signIn.Range("A13").Select()
